# Weekly data update: insert a new price record as row 374 in the
# "Terminal La Palmera de La Serena - Zanahoria" sheet, pushing the
# existing rows 374:490 down to 375:491.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 374; this shifts rows 374:490 down
# to 375:491 and extends the sheet's used range to A1:R491.
$ws.Rows("374").Insert()

# Populate the newly inserted row 374 with the new weekly record.
$ws.Range("A374").Value = 8
$ws.Range("B374").Value = "Terminal La Palmera de La Serena"
$ws.Range("C374").Value = "Coquimbo"
$ws.Range("D374").Value = 44985
$ws.Range("E374").Value = 4
$ws.Range("F374").Value = 100114013
$ws.Range("G374").Value = "Zanahoria"
$ws.Range("H374").Value = "Sin especificar"
$ws.Range("I374").Value = "Primera"
$ws.Range("J374").Value = 600
$ws.Range("K374").Value = 6000
$ws.Range("L374").Value = 7000
$ws.Range("M374").Value = 6500
$ws.Range("N374").Value = "$/saco 20 kilos"
$ws.Range("O374").Value = "Provincia del Elquí"
$ws.Range("P374").Value = 325
$ws.Range("Q374").Value = 20
$ws.Range("R374").Value = "Hortaliza"
